$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a rolling 28/29-day window of per-symbol quantities, one
# row per date (column A) with computed columns B:J. This update rolls the
# window forward: the 3 oldest dates (2024-09-10..12) drop off, and 2 new
# dates (2024-10-09, 2024-10-10) are appended, so every surviving row's
# data shifts up by three rows and two fresh rows are appended at the
# bottom (net: one fewer data row, 29 -> 28).

# 1) Shift rows 5..30 (old) up to rows 2..27 (new) in one shot, carrying
#    values *and* formatting (so column A keeps its shared-string/style
#    bits exactly as Excel would when you drag a block up).
$ws.Range("A5:J30").Copy($ws.Range("A2:J27"))

# 2) Build the two new trailing rows (28: 2024-10-09, 29: 2024-10-10).
#    They repeat the same B:J figures as the new last "old" row (27).
$ws.Range("A28:J28").Value = $ws.Range("A27:J27").Value2
$ws.Range("A29:J29").Value = $ws.Range("A27:J27").Value2

# Column A must hold literal date-text (matching the sheet's existing
# "yyyy-mm-dd" strings) rather than an auto-converted date serial, so we
# stage it via a formula (never auto-converted) in a scratch cell, then
# paste just the computed value back in - this keeps the destination's
# existing style/shared-string typing intact.
$ws.Range("L1").Formula = '="2024-10-09"'
$ws.Range("L1").Copy()
$ws.Range("A28").PasteSpecial(-4163)

$ws.Range("L1").Formula = '="2024-10-10"'
$ws.Range("L1").Copy()
$ws.Range("A29").PasteSpecial(-4163)

# Clean up the scratch cell and the now-stale last row of the old window.
$ws.Range("L1").ClearContents()
$ws.Range("A30:J30").ClearContents()

$excel.CutCopyMode = $false
